$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main root")
$ws.Range("A2").Select()
